# Update public EPEX spot prices workbook
#  - "Prix Spot" sheet: insert a new day column (24-nov) before the
#    01-oct. column, shifting every subsequent date column one place
#    to the right.
#  - "Gaz" and "CO2" sheets: append the latest two daily quotes.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "Prix Spot" sheet — insert new "24-nov" column before "01-oct."
# ---------------------------------------------------------------------
$wsPrix = $wb.Worksheets.Item("Prix Spot")

# Column DY currently holds "01-oct." (and the daily values below it).
# Inserting a whole column here pushes DY:FC -> DZ:FD and leaves a
# blank DY column ready to receive the new day's data.
$wsPrix.Range("DY1").EntireColumn.Insert()

$wsPrix.Range("DY1").Value = "24-nov"
$wsPrix.Range("DY2:DY25").Value = "-"

# ---------------------------------------------------------------------
# 2) "Gaz" sheet — append 2025-11-22 and 2025-11-23
# ---------------------------------------------------------------------
$wsGaz = $wb.Worksheets.Item("Gaz")

$wsGaz.Range("A158").NumberFormat = "@"
$wsGaz.Range("A158").Value = "2025-11-22"
$wsGaz.Range("A158").ClearFormats()
$wsGaz.Range("B158").Value = 29.2

$wsGaz.Range("A159").NumberFormat = "@"
$wsGaz.Range("A159").Value = "2025-11-23"
$wsGaz.Range("A159").ClearFormats()
$wsGaz.Range("B159").Value = 29.2

# ---------------------------------------------------------------------
# 3) "CO2" sheet — append 2025-11-22 and 2025-11-23
# ---------------------------------------------------------------------
$wsCO2 = $wb.Worksheets.Item("CO2")

$wsCO2.Range("A158").NumberFormat = "@"
$wsCO2.Range("A158").Value = "2025-11-22"
$wsCO2.Range("A158").ClearFormats()
$wsCO2.Range("B158").Value = 80.28

$wsCO2.Range("A159").NumberFormat = "@"
$wsCO2.Range("A159").Value = "2025-11-23"
$wsCO2.Range("A159").ClearFormats()
$wsCO2.Range("B159").Value = 80.28
